# Updated symbol list (cryptos.xlsx) - refresh Price (D) and Volume(1h) (E)
# columns with the latest scraped quotes, GitHub Actions style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns are stored as plain text (e.g. "290.83", "-6.29%")
# rather than numbers, so force text storage before writing the new quotes.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "290.83"
$ws.Range("E2").Value = "-6.29%"
$ws.Range("D3").Value = "39.52"
$ws.Range("E3").Value = "-3.86%"
$ws.Range("D4").Value = "5.015"
$ws.Range("E4").Value = "-3.85%"
$ws.Range("D5").Value = "0.07352"
$ws.Range("E5").Value = "-4.32%"
$ws.Range("D6").Value = "4.287"
$ws.Range("E6").Value = "-0.25%"
$ws.Range("D7").Value = "1.555"
$ws.Range("E7").Value = "-10.93%"
$ws.Range("D8").Value = "0.9139"
$ws.Range("E8").Value = "-1.19%"
$ws.Range("D9").Value = "0.1189"
$ws.Range("E9").Value = "-6.40%"
$ws.Range("D10").Value = "0.1731"
$ws.Range("E10").Value = "-5.04%"
$ws.Range("D11").Value = "0.08707"
$ws.Range("E11").Value = "-4.38%"
$ws.Range("D12").Value = "0.04162"
$ws.Range("E12").Value = "-0.12%"
$ws.Range("D13").Value = "0.1051"
$ws.Range("E13").Value = "0.00%"
$ws.Range("D14").Value = "0.001276"
$ws.Range("E14").Value = "-1.17%"
$ws.Range("D15").Value = "0.005855"
$ws.Range("E15").Value = "-0.57%"
$ws.Range("D16").Value = "3.383"
$ws.Range("E16").Value = "0.87%"
$ws.Range("D18").Value = "0.3297"
$ws.Range("E18").Value = "-0.69%"
$ws.Range("D19").Value = "7.535"
$ws.Range("E19").Value = "2.02%"
$ws.Range("D20").Value = "0.1354"
$ws.Range("E20").Value = "0.03%"
$ws.Range("D21").Value = "0.2884"
$ws.Range("E21").Value = "5.90%"
$ws.Range("D22").Value = "0.03843"
$ws.Range("E22").Value = "-4.20%"
$ws.Range("E23").Value = "0.43%"
$ws.Range("D24").Value = "0.003683"
$ws.Range("E24").Value = "-10.08%"
$ws.Range("D25").Value = "0.0001282"
$ws.Range("E25").Value = "0.96%"
$ws.Range("D26").Value = "0.0003727"
$ws.Range("E26").Value = "-95.04%"
$ws.Range("D38").Value = "0.02310"
$ws.Range("E38").Value = "-8.58%"
$ws.Range("D39").Value = "0.05011"
$ws.Range("E39").Value = "-5.66%"
$ws.Range("D40").Value = "0.007705"
$ws.Range("E40").Value = "-1.81%"
$ws.Range("E41").Value = "149.34%"
$ws.Range("D42").Value = "0.1268"
$ws.Range("E42").Value = "-3.11%"
$ws.Range("D43").Value = "0.007379"
$ws.Range("E43").Value = "11.21%"
$ws.Range("D44").Value = "0.007673"
$ws.Range("E44").Value = "-5.51%"
$ws.Range("D45").Value = "0.3162"
$ws.Range("E45").Value = "2.70%"
$ws.Range("D46").Value = "0.00006545"
$ws.Range("E46").Value = "-3.45%"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").Value = "0.01%"
$ws.Range("E48").Value = "10.83%"
$ws.Range("D49").Value = "0.004206"
$ws.Range("E49").Value = "35.51%"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").Value = "0.01%"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").Value = "0.01%"

# Restore default (no explicit) cell format so styling matches the source file
$ws.Range("D2:E51").ClearFormats()
